$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$E2Value = @'
[Simon D%Harding%NULL%1,            Joanna L%Sharman%NULL%2,            Joanna L%Sharman%NULL%0,            Elena%Faccenda%NULL%1,            Chris%Southan%NULL%1,            Adam J%Pawson%NULL%1,            Sam%Ireland%NULL%1,            Alasdair J G%Gray%NULL%1,            Liam%Bruce%NULL%1,            Stephen P H%Alexander%NULL%1,            Stephen%Anderton%NULL%1,            Clare%Bryant%NULL%1,            Anthony P%Davenport%NULL%1,            Christian%Doerig%NULL%1,            Doriano%Fabbro%NULL%1,            Francesca%Levi-Schaffer%NULL%1,            Michael%Spedding%NULL%1,            Jamie A%Davies%jamie.davies@ed.ac.uk%1,            NULL%NULL%NULL%0]
'@
$ws.Range("E2").Value = $E2Value

$E3Value = @'
[ Ourania S.%Kotsiou%null%1,             Sotirios G.%Zarogiannis%null%1,             Konstantinos I.%Gourgoulianis%null%1,           Ourania S.%Kotsiou%null%1,           Sotirios G.%Zarogiannis%null%1,           Konstantinos I.%Gourgoulianis%null%1]
'@
$ws.Range("E3").Value = $E3Value

$E4Value = @'
[ Jonathan%Messika%null%1,             Benjamin%Sztrymf%null%1,             Fabrice%Bertrand%null%1,             Typhaine%Billard-Pomares%null%1,             Guilene%Barnaud%null%1,             Catherine%Branger%null%1,             Didier%Dreyfuss%null%1,             Jean-Damien%Ricard%null%1,           Jonathan%Messika%null%1,           Benjamin%Sztrymf%null%1,           Fabrice%Bertrand%null%1,           Typhaine%Billard-Pomares%null%1,           Guilene%Barnaud%null%1,           Catherine%Branger%null%1,           Didier%Dreyfuss%null%1,           Jean-Damien%Ricard%null%1]
'@
$ws.Range("E4").Value = $E4Value

$E5Value = @'
[Carl%Llor%carles.llor@urv.cat%1,            Ana%Moragas%amoragasm@meditex.es%1,            Carolina%Bayona%cbayona.tarte.ics@gencat.cat%1,            Rosa%Morros%rmorros@idiapjgol.org%1,            Helena%Pera%hpera@idiapjgol.org%1,            Josep M%Cots%23465jcy@comb.cat%1,            Yvonne%Fernández%yfernandez.tarte.ics@gencat.cat%1,            Marc%Miravitlles%marcm@separ.es%1,            Albert%Boada%aboadav.bcn.ics@gencat.cat%1]
'@
$ws.Range("E5").Value = $E5Value

$E6Value = @'
[Holly%Epperly%xref no email%1,    Frances L.%Vaughn%xref no email%1,    Andrew D.%Mosholder%xref no email%1,    Elizabeth M.%Maloney%xref no email%1,    Lewis%Rubinson%xref no email%1]
'@
$ws.Range("E6").Value = $E6Value

$E7Value = @'
[ G.%Voiriot%null%1,             M.%Chalumeau%null%1,             J.%Messika%null%1,             D.%Basille%null%1,             B.%Philippe%null%1,             J.-D.%Ricard%null%1,             C.%Andrejak%null%1,             V.%Jounieaux%null%1,             O.%Sanchez%null%1,             M.%Fartoukh%null%1,           G.%Voiriot%null%1,           M.%Chalumeau%null%1,           J.%Messika%null%1,           D.%Basille%null%1,           B.%Philippe%null%1,           J.-D.%Ricard%null%1,           C.%Andrejak%null%1,           V.%Jounieaux%null%1,           O.%Sanchez%null%1,           M.%Fartoukh%null%1]
'@
$ws.Range("E7").Value = $E7Value

$D8Value = @'

Purpose
Non-steroidal anti-inflammatory drugs (NSAIDs) are frequently prescribed or used as self-medication in cases of community-acquired pneumonia (CAP).
 Nevertheless, the consequences of such medication on the risk of pleuroparenchymal complications are not well known.
 The aim was to investigate whether exposure to NSAIDs prior to hospital admission among patients suffering from CAP is associated with the development of pleural complications or a lung abscess.
Methods
All consecutive non-immunocompromised patients with CAP and admitted to a university hospital were prospectively included (2-year period).
 The risk of pleuropulmonary complications was analyzed according to previous exposure to NSAIDs.
Results
Of the 221 included patients, 40 (18.1%) had developed a pleuropulmonary complication.
 NSAIDs intake prior to admission was reported for 24 patients (10.9%) who were younger (50.6 ± 18.5 vs.
 66.5 ± 16.4 years; p = 0.001), had less comorbidities (60 vs.
 25.1%; p = 0.001), had a longer duration between the first symptoms of CAP and the start of an antibiotic therapy (6.1 ± 7.6 vs.
 2.8 ± 3.8 days; p = 0.001), and who had a higher incidence of pleuropulmonary complications (33.3 vs.
 16.2%; p = 0.048).
 In multivariate analyses, two factors were independently associated with the development of pleuroparenchymal complications: NSAIDs intake [Odds Ratio (OR) = 2.57 [1.02–6.64]; p = 0.049] and alcohol abuse (OR = 2.68 [1.27–5.69]; p = 0.01).
Conclusions
Our findings suggest that NSAIDs, often taken by young and healthy patients, may worsen the course of CAP with delayed therapy and a higher rate of pleuropulmonary complications.

'@
$ws.Range("D8").Value = $D8Value

$E8Value = @'
[ Damien%Basille%null%0,  Nathalie%Plouvier%null%1,  Charlotte%Trouve%null%1,  Pierre%Duhaut%null%1,  Claire%Andrejak%null%1,  Vincent%Jounieaux%null%1]
'@
$ws.Range("E8").Value = $E8Value

$E9Value = @'
[Damien%Basille%xref no email%1,    Reimar Wernich%Thomsen%xref no email%1,    Morten%Madsen%xref no email%1,    Pierre%Duhaut%xref no email%1,    Claire%Andrejak%xref no email%1,    Vincent%Jounieaux%xref no email%1,    Henrik Toft%S\u00f8rensen%xref no email%1]
'@
$ws.Range("E9").Value = $E9Value
